# Scheduled data refresh: update market price / profit figures across item sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 621.8095
$ws.Range("I127").Value = 324.25
$ws.Range("J127").Value = 1018.55554
$ws.Range("K127").Value = 972.75
$ws.Range("L127").Value = 3055.66662
$ws.Range("M127").Value = 3987.25
$ws.Range("N127").Value = -12975.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1255.7407
$ws.Range("I45").Value = 1026.579
$ws.Range("K45").Value = 1026.579
$ws.Range("M45").Value = -649.579

$ws.Range("H97").Value = 878.67566
$ws.Range("I97").Value = 878.67566
$ws.Range("K97").Value = 878.67566
$ws.Range("M97").Value = -382.67566

$ws.Range("H122").Value = 4293.857
$ws.Range("I122").Value = 4691.75
$ws.Range("J122").Value = 3763.3333
$ws.Range("K122").Value = 14075.25
$ws.Range("L122").Value = 11289.9999
$ws.Range("M122").Value = -11625.25
$ws.Range("N122").Value = -16189.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1803.7805
$ws.Range("I105").Value = 1709.4445
$ws.Range("J105").Value = 1985.7142
$ws.Range("K105").Value = 1709.4445
$ws.Range("L105").Value = 1985.7142
$ws.Range("M105").Value = 37.55549999999994
$ws.Range("N105").Value = -5479.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7376.846
$ws.Range("I86").Value = 4649.8335
$ws.Range("J86").Value = 9714.286
$ws.Range("K86").Value = 4649.8335
$ws.Range("L86").Value = 9714.286
$ws.Range("M86").Value = -3526.8335
$ws.Range("N86").Value = -11960.286

$ws.Range("H89").Value = 7376.846
$ws.Range("I89").Value = 4649.8335
$ws.Range("J89").Value = 9714.286
$ws.Range("K89").Value = 23249.1675
$ws.Range("L89").Value = 48571.43
$ws.Range("M89").Value = -17633.1675
$ws.Range("N89").Value = -59803.43

$ws.Range("H107").Value = 232.42857
$ws.Range("I107").Value = 180
$ws.Range("K107").Value = 180
$ws.Range("M107").Value = 1740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1571.8182
$ws.Range("I139").Value = 748.75
$ws.Range("K139").Value = 2246.25
$ws.Range("M139").Value = 2893.75

$ws.Range("H140").Value = 1481.9
$ws.Range("I140").Value = 959.8570999999999
$ws.Range("J140").Value = 2700
$ws.Range("K140").Value = 2879.5713
$ws.Range("L140").Value = 8100
$ws.Range("M140").Value = 2300.4287
$ws.Range("N140").Value = -18460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5415.4287
$ws.Range("I70").Value = 3965.849
$ws.Range("J70").Value = 8616.583000000001
$ws.Range("K70").Value = 3965.849
$ws.Range("L70").Value = 8616.583000000001
$ws.Range("M70").Value = -3695.849
$ws.Range("N70").Value = -9156.583000000001

$ws.Range("H73").Value = 5415.4287
$ws.Range("I73").Value = 3965.849
$ws.Range("J73").Value = 8616.583000000001
$ws.Range("K73").Value = 3965.849
$ws.Range("L73").Value = 8616.583000000001
$ws.Range("M73").Value = -3029.849
$ws.Range("N73").Value = -10488.583

$ws.Range("H97").Value = 1036.8
$ws.Range("I97").Value = 921
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 921
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -425
$ws.Range("N97").Value = -2492

$ws.Range("H122").Value = 1257.8
$ws.Range("I122").Value = 1114
$ws.Range("J122").Value = 1593.3334
$ws.Range("K122").Value = 3342
$ws.Range("L122").Value = 4780.0002
$ws.Range("M122").Value = -892
$ws.Range("N122").Value = -9680.0002

$ws.Range("H123").Value = 31124.75
$ws.Range("J123").Value = 31124.75
$ws.Range("L123").Value = 31124.75
$ws.Range("N123").Value = -36024.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 9636.333000000001
$ws.Range("I4").Value = 8909
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 8909
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -8796
$ws.Range("N4").Value = -10226

$ws.Range("H28").Value = 9636.333000000001
$ws.Range("I28").Value = 8909
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 8909
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -8677
$ws.Range("N28").Value = -10464

$ws.Range("H37").Value = 9636.333000000001
$ws.Range("I37").Value = 8909
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 8909
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -8802
$ws.Range("N37").Value = -10214

$ws.Range("H68").Value = 2279.5833
$ws.Range("I68").Value = 1758
$ws.Range("J68").Value = 2652.1428
$ws.Range("K68").Value = 1758
$ws.Range("L68").Value = 2652.1428
$ws.Range("M68").Value = -1009
$ws.Range("N68").Value = -4150.1428

$ws.Range("H71").Value = 2279.5833
$ws.Range("I71").Value = 1758
$ws.Range("J71").Value = 2652.1428
$ws.Range("K71").Value = 8790
$ws.Range("L71").Value = 13260.714
$ws.Range("M71").Value = -5046
$ws.Range("N71").Value = -20748.714

$ws.Range("H100").Value = 40003508
$ws.Range("I100").Value = 4506.9375
$ws.Range("J100").Value = 111112850
$ws.Range("K100").Value = 4506.9375
$ws.Range("L100").Value = 111112850
$ws.Range("M100").Value = -3965.9375
$ws.Range("N100").Value = -111113932

$ws.Range("H122").Value = 3089.2942
$ws.Range("I122").Value = 2578.5715
$ws.Range("J122").Value = 3446.8
$ws.Range("K122").Value = 7735.7145
$ws.Range("L122").Value = 10340.4
$ws.Range("M122").Value = -5285.7145
$ws.Range("N122").Value = -15240.4

$ws.Range("H133").Value = 22738
$ws.Range("J133").Value = 22738
$ws.Range("L133").Value = 22738
$ws.Range("N133").Value = -27798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 8338938.5
$ws.Range("I21").Value = 16669393
$ws.Range("J21").Value = 8483.333000000001
$ws.Range("K21").Value = 16669393
$ws.Range("L21").Value = 8483.333000000001
$ws.Range("M21").Value = -16669158
$ws.Range("N21").Value = -8953.333000000001

$ws.Range("H35").Value = 8338938.5
$ws.Range("I35").Value = 16669393
$ws.Range("J35").Value = 8483.333000000001
$ws.Range("K35").Value = 16669393
$ws.Range("L35").Value = 8483.333000000001
$ws.Range("M35").Value = -16669103
$ws.Range("N35").Value = -9063.333000000001

$ws.Range("H113").Value = 369.32144
$ws.Range("I113").Value = 337.05
$ws.Range("J113").Value = 450
$ws.Range("K113").Value = 1011.15
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 1158.85
$ws.Range("N113").Value = -5690

$ws.Range("H122").Value = 2199.1428
$ws.Range("I122").Value = 2315.6667
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6947.000100000001
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4497.000100000001
$ws.Range("N122").Value = -9400

$ws.Range("H126").Value = 962.0833
$ws.Range("I126").Value = 742.6786
$ws.Range("J126").Value = 1730
$ws.Range("K126").Value = 2228.0358
$ws.Range("L126").Value = 5190
$ws.Range("M126").Value = 241.9642000000003
$ws.Range("N126").Value = -10130

$ws.Range("H135").Value = 29053
$ws.Range("J135").Value = 29053
$ws.Range("L135").Value = 29053
$ws.Range("N135").Value = -39193
